$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of gradient data to append (distance, MEAN, STD, MIN, MAX, COUNT, Month)
$newRows = @(
    @(1100, 42.23291778564453, 3.433120489120483, 26.42891883850098, 52.21511077880859, 18476, "06"),
    @(1200, 42.18862915039062, 3.544787168502808, 26.11795234680176, 52.32616806030273, 18568, "06"),
    @(1300, 42.10209274291992, 3.621676445007324, 25.38154602050781, 52.63200759887695, 18438, "06"),
    @(1400, 41.99977111816406, 3.662901878356934, 25.37300300598145, 53.89978790283203, 18366, "06"),
    @(1500, 41.92996978759766, 3.78806734085083, 25.78819274902344, 56.7018928527832, 18392, "06")
)

$startRow = 21
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $gCell = $ws.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $data[6]
}
